$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently holds Julia's data; replace it with what used to be row 4 (Brendon/admin)
$ws.Range("A3").Value = "Brendon"
$ws.Range("B3").Value = "$2b$12$.BvYnlk164DiJ9jTsNHK7OAkq045fG9Ma6Vk/mWfS6pBozLIyjlMS"
$ws.Range("C3").Value = "admin"

# Delete the now-duplicate row 4 entirely, shifting rows up
$ws.Rows("4:4").Delete()
